# Update line-states / contingency table (case 1_36) with the "rene fine"
# re-run values, and append two new contingency rows (extr7, extr8).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Existing rows: refreshed C/D/E values -------------------------------

# row 8  (name=line5)
$ws.Range("C8").Value = 14
$ws.Range("D8").Value = 11

# row 9  (name=line6)
$ws.Range("C9").Value = 16

# row 10 (name=extr1)
$ws.Range("C10").Value = 5
$ws.Range("D10").Value = 12
$ws.Range("E10").Value = $true

# row 11 (name=extr2)
$ws.Range("C11").Value = 5
$ws.Range("D11").Value = 9

# row 12 (name=extr3)
$ws.Range("C12").Value = 10

# row 13 (name=extr4)
$ws.Range("D13").Value = 8
$ws.Range("E13").Value = $true

# row 14 (name=extr5)
$ws.Range("C14").Value = 9
$ws.Range("D14").Value = 11

# row 15 (name=extr6)
$ws.Range("C15").Value = 7
$ws.Range("D15").Value = 11
$ws.Range("E15").Value = $false

# --- New rows 16 & 17 (extr7, extr8) --------------------------------------

$ws.Range("A16").Value = 14
$ws.Range("B16").Value = "extr7"
$ws.Range("C16").Value = 5
$ws.Range("D16").Value = 7
$ws.Range("E16").Value = $true

$ws.Range("A17").Value = 15
$ws.Range("B17").Value = "extr8"
$ws.Range("C17").Value = 8
$ws.Range("D17").Value = 5
$ws.Range("E17").Value = $true

# Match the formatting (bold, centered, bordered) already used by column A
# in the earlier data rows.
$ws.Range("A15").Copy()
$ws.Range("A16:A17").PasteSpecial(-4122)
$excel.CutCopyMode = $false
